# Remove the "state" (country) column from the "all" sheet, shifting the
# remaining atco_ops / support / all columns left by one.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("all")
$ws.Columns.Item(2).Delete()

# Mirror the author's final selection (whole of the new column B selected).
$ws.Range("B1:B1048576").Select()
